# Update crypto price/volume data per daily GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.343.27'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.928.58'
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.10'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.45'
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("D9").Value = '2.928.37'
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.88'
$ws.Range("E10").Value = '  +1.97%  '
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.444'
$ws.Range("E12").Value = '  -1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000228'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.88'
$ws.Range("E14").Value = '  -1.96%  '
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '3.414.03'
$ws.Range("E16").Value = '  +0.15%  '
$ws.Range("D17").Value = '61.228.98'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.73'
$ws.Range("E18").Value = '  -1.64%  '
$ws.Range("D19").Value = '2.928.15'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '432.43'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("E22").Value = '  +1.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.11'
$ws.Range("E23").Value = '  -1.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '80.95'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("B25").Value = 'RenderToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.96'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("E26").Value = '  +2.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.12'
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.35'
$ws.Range("E29").Value = '  +8.14%  '
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.63'
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.21'
$ws.Range("E32").Value = '  -1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.62'
$ws.Range("E33").Value = '  -0.74%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  +1.71%  '
$ws.Range("D35").Value = '0.0₃0872'
$ws.Range("E35").Value = '  +3.41%  '
$ws.Range("E36").Value = '  +0.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.13'
$ws.Range("E37").Value = '  +3.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.65'
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '50.01'
$ws.Range("E39").Value = '  +0.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.02'
$ws.Range("E40").Value = '  -0.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.125'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.63'
$ws.Range("E42").Value = '  -1.87%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.291'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.73'
$ws.Range("E44").Value = '  -5.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '381.60'
$ws.Range("E45").Value = '  +1.57%  '
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("D47").Value = '2.717.17'
$ws.Range("E47").Value = '  +2.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.70'
$ws.Range("E48").Value = '  -1.69%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.34'
$ws.Range("E50").Value = '  -4.86%  '
$ws.Range("E51").Value = '  +0.14%  '
